# Atualizado por script em 21-12-2023 20:46
#
# 1) Three pairs of adjacent rows had their match-detail columns (F:V)
#    swapped while the index/metadata columns (A:E) stayed put.
# 2) One brand-new match row (128) was appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowDetails {
    param($rowA, $rowB)

    $rangeA = $ws.Range("F$rowA`:V$rowA")
    $rangeB = $ws.Range("F$rowB`:V$rowB")

    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

# --- Swap the three pairs of mismatched rows -------------------------------
Swap-RowDetails 32 33
Swap-RowDetails 60 61
Swap-RowDetails 111 112

# --- Append the new match row 128 ------------------------------------------
$newRow = 128
$prevRow = 127

$ws.Range("A$newRow").Value = 127
$ws.Range("B$newRow").Value = "portugal"
$ws.Range("C$newRow").Value = "liga-portugal-2"
$ws.Range("D$newRow").Value = "2023-2024"
$ws.Range("E$newRow").Value = 45281.79166666666
$ws.Range("F$newRow").Value = "Torreense"
$ws.Range("G$newRow").Value = 3
$ws.Range("H$newRow").Value = "Vilaverdense"
$ws.Range("I$newRow").Value = 1
$ws.Range("J$newRow").Value = 1.72
$ws.Range("K$newRow").Value = "17/12/2023 12:12"
$ws.Range("L$newRow").Value = 1.71
$ws.Range("M$newRow").Value = "21/12/2023 18:56"
$ws.Range("N$newRow").Value = 3.67
$ws.Range("O$newRow").Value = "17/12/2023 12:12"
$ws.Range("P$newRow").Value = 3.67
$ws.Range("Q$newRow").Value = "21/12/2023 18:56"
$ws.Range("R$newRow").Value = 4.66
$ws.Range("S$newRow").Value = "17/12/2023 12:12"
$ws.Range("T$newRow").Value = 5.57
$ws.Range("U$newRow").Value = "21/12/2023 18:56"
$ws.Range("V$newRow").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/torreense-vilaverdense-fc/Cb2nDCv0/"

# Carry over the same per-column styling used by every other data row
# (bold/centered/bordered index in column A, datetime format in column E)
# by copying the formatting-only from the row directly above.
$ws.Range("A$prevRow").Copy()
$ws.Range("A$newRow").PasteSpecial(-4122)

$ws.Range("E$prevRow").Copy()
$ws.Range("E$newRow").PasteSpecial(-4122)

$excel.CutCopyMode = 0
